$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Content change: "No Subnet Data" -> "No Subnet" in C1 ---
$ws.Range("C1").Value = "No Subnet"

# --- Formatting artifacts left over on C1:F1 (red font highlight, with E1
#     reverted back to the normal/black font) ---
$ws.Range("C1:D1").Font.Color = 255
$ws.Range("E1").Font.Color = 0
$ws.Range("F1").Font.Color = 255

# --- Page setup: paper size / orientation now explicit on the sheet ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# --- View state: scrolled over one column, new active cell/selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("D18").Select() | Out-Null
